$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.943.77'
$ws.Range("D3").Value = '1.633.20'
$ws.Range("E3").Value = '  +1.75%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.83%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '28.68'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.26%  '
$ws.Range("E9").Value = '  +2.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0609'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.98%  '
$ws.Range("E11").Value = '  +0.21%  '
$ws.Range("D12").Value = '1.866.97'
$ws.Range("E12").Value = '  +1.74%  '
$ws.Range("D13").Value = '1.632.35'
$ws.Range("E13").Value = '  +1.71%  '
$ws.Range("E14").Value = '  +2.94%  '
$ws.Range("E15").Value = '  +18.21%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.86'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.64%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '29.937.32'
$ws.Range("E17").Value = '  +0.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.12'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.94%  '
$ws.Range("D20").Value = '0.0₃0702'
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.85'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.68%  '
$ws.Range("E23").Value = '  +2.51%  '
$ws.Range("E24").Value = '  +1.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.73'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.53'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.48%  '
$ws.Range("E27").Value = '  +1.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.60'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0486'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.11'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.21%  '
$ws.Range("E32").Value = '  +4.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.17'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.42%  '
$ws.Range("D34").Value = '1.424.28'
$ws.Range("E34").Value = '  -0.11%  '
$ws.Range("E35").Value = '  +4.61%  '
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.43%  '
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("E39").Value = '  +0.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '75.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +14.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.552'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.78%  '
$ws.Range("E42").Value = '  +2.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.829'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0489'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '53.14'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.14%  '
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("E47").Value = '  +3.24%  '
$ws.Range("D48").Value = '1.775.23'
$ws.Range("E48").Value = '  +1.97%  '
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₆0112'
$ws.Range("E50").Value = '  +7.90%  '
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '89.31'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.04%  '
